# Poland Ekstraklasa update (31-03-2024 20:29) -------------------------------
# Old rows 226-230 (match ids 224-228, i.e. payload for B..AC) are pushed
# down by 4 rows to become rows 230-234; four brand-new match rows are
# written in the freed rows 226-229. Column A (the sequential "id" index,
# always row-2) never moves and is already correct for rows 226-230, so we
# only need to add it for the brand new rows 231-234.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 226-230: overwrite the data payload (B..AC) in place.
#    Rows 226-229 additionally gain H/I/J (full-time score/result) and
#    AB/AC (closing Asian-handicap over/under) which did not exist before.
#    Row 230 keeps the exact same (no H/I/J/AB/AC) shape the old row 226 had.
# ---------------------------------------------------------------------------

# Row 226 (new match 6774469)
$ws.Range("B226").Value2 = 6774469
$ws.Range("E226").Value2 = 45381.35416666666
$ws.Range("F226").Value = "Rakow Czestochowa"
$ws.Range("G226").Value = "Ruch Chorzow"
$ws.Range("H226").Value2 = 1
$ws.Range("I226").Value2 = 1
$ws.Range("J226").Value = "D"
$ws.Range("K226").Value2 = 1.4
$ws.Range("L226").Value2 = 4.75
$ws.Range("M226").Value2 = 7.5
$ws.Range("N226").Value2 = 1.4
$ws.Range("O226").Value2 = 4.75
$ws.Range("P226").Value2 = 7
$ws.Range("Q226").Value2 = -1.25
$ws.Range("R226").Value2 = 1.975
$ws.Range("S226").Value2 = 1.875
$ws.Range("T226").Value2 = 2.75
$ws.Range("U226").Value2 = 1.9
$ws.Range("V226").Value2 = 1.95
$ws.Range("W226").Value2 = -1
$ws.Range("X226").Value2 = 3.75
$ws.Range("Y226").Value2 = -1
$ws.Range("Z226").Value2 = -1
$ws.Range("AA226").Value2 = 0.875
$ws.Range("AB226").Value2 = -1
$ws.Range("AC226").Value2 = 0.95

# Row 227 (new match 6774468)
$ws.Range("B227").Value2 = 6774468
$ws.Range("E227").Value2 = 45381.45833333334
$ws.Range("F227").Value = "Jagiellonia Bialystok"
$ws.Range("G227").Value = "LKS Lodz"
$ws.Range("H227").Value2 = 6
$ws.Range("I227").Value2 = 0
$ws.Range("J227").Value = "H"
$ws.Range("K227").Value2 = 1.4
$ws.Range("L227").Value2 = 5
$ws.Range("M227").Value2 = 7
$ws.Range("N227").Value2 = 1.333
$ws.Range("O227").Value2 = 5.75
$ws.Range("P227").Value2 = 8
$ws.Range("Q227").Value2 = -1.5
$ws.Range("R227").Value2 = 1.85
$ws.Range("S227").Value2 = 2
$ws.Range("T227").Value2 = 3.25
$ws.Range("U227").Value2 = 2.025
$ws.Range("V227").Value2 = 1.825
$ws.Range("W227").Value2 = 0.333
$ws.Range("X227").Value2 = -1
$ws.Range("Y227").Value2 = -1
$ws.Range("Z227").Value2 = 0.8500000000000001
$ws.Range("AA227").Value2 = -1
$ws.Range("AB227").Value2 = 1.025
$ws.Range("AC227").Value2 = -1

# Row 228 (new match 6775574)
$ws.Range("B228").Value2 = 6775574
$ws.Range("E228").Value2 = 45381.5625
$ws.Range("F228").Value = "Piast Gliwice"
$ws.Range("G228").Value = "Slask Wroclaw"
$ws.Range("H228").Value2 = 2
$ws.Range("I228").Value2 = 2
$ws.Range("J228").Value = "D"
$ws.Range("K228").Value2 = 2.1
$ws.Range("L228").Value2 = 3.1
$ws.Range("M228").Value2 = 4
$ws.Range("N228").Value2 = 2.2
$ws.Range("O228").Value2 = 2.9
$ws.Range("P228").Value2 = 4
$ws.Range("Q228").Value2 = -0.25
$ws.Range("R228").Value2 = 1.825
$ws.Range("S228").Value2 = 2.025
$ws.Range("T228").Value2 = 1.75
$ws.Range("U228").Value2 = 1.875
$ws.Range("V228").Value2 = 1.975
$ws.Range("W228").Value2 = -1
$ws.Range("X228").Value2 = 1.9
$ws.Range("Y228").Value2 = -1
$ws.Range("Z228").Value2 = -0.5
$ws.Range("AA228").Value2 = 0.5125
$ws.Range("AB228").Value2 = 0.875
$ws.Range("AC228").Value2 = -1

# Row 229 (new match 6775575)
$ws.Range("B229").Value2 = 6775575
$ws.Range("E229").Value2 = 45381.66666666666
$ws.Range("F229").Value = "Pogon Szczecin"
$ws.Range("G229").Value = "Cracovia Krakow"
$ws.Range("H229").Value2 = 3
$ws.Range("I229").Value2 = 1
$ws.Range("J229").Value = "H"
$ws.Range("K229").Value2 = 1.909
$ws.Range("L229").Value2 = 3.6
$ws.Range("M229").Value2 = 4
$ws.Range("N229").Value2 = 1.727
$ws.Range("O229").Value2 = 3.75
$ws.Range("P229").Value2 = 5
$ws.Range("Q229").Value2 = -0.75
$ws.Range("R229").Value2 = 1.95
$ws.Range("S229").Value2 = 1.9
$ws.Range("T229").Value2 = 2.5
$ws.Range("U229").Value2 = 1.875
$ws.Range("V229").Value2 = 1.975
$ws.Range("W229").Value2 = 0.7270000000000001
$ws.Range("X229").Value2 = -1
$ws.Range("Y229").Value2 = -1
$ws.Range("Z229").Value2 = 0.95
$ws.Range("AA229").Value2 = -1
$ws.Range("AB229").Value2 = 0.875
$ws.Range("AC229").Value2 = -1

# Row 230 (now holds what used to be row 226's match, 6774877)
$ws.Range("B230").Value2 = 6774877
$ws.Range("E230").Value2 = 45383.3125
$ws.Range("F230").Value = "Puszcza Niepolomice"
$ws.Range("G230").Value = "Radomiak Radom"
$ws.Range("K230").Value2 = 2.625
$ws.Range("L230").Value2 = 3.4
$ws.Range("M230").Value2 = 2.6
$ws.Range("N230").Value2 = 2.625
$ws.Range("O230").Value2 = 3.4
$ws.Range("P230").Value2 = 2.6
$ws.Range("Q230").Value2 = 0
$ws.Range("R230").Value2 = 1.975
$ws.Range("S230").Value2 = 1.875
$ws.Range("T230").Value2 = 2.25
$ws.Range("U230").Value2 = 1.825
$ws.Range("V230").Value2 = 2.025
$ws.Range("W230").Value2 = 0
$ws.Range("X230").Value2 = 0
$ws.Range("Y230").Value2 = 0
$ws.Range("Z230").Value2 = 0
$ws.Range("AA230").Value2 = 0

# ---------------------------------------------------------------------------
# 2) Rows 231-234 are brand new - write the full row (A..AC, no H/I/J/AB/AC,
#    matching the shape the old rows 227-230 had) and copy the existing
#    formatting for the A (bold+bordered id) and E (date) columns down from
#    row 230 so no new style entries are introduced.
# ---------------------------------------------------------------------------

$ws.Range("A230:E230").Copy()
$ws.Range("A231:E234").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 231 (was old row 227, match 6775576)
$ws.Range("A231").Value2 = 229
$ws.Range("B231").Value2 = 6775576
$ws.Range("C231").Value = "Poland Ekstraklasa"
$ws.Range("D231").Value = "Poland Ekstraklasa"
$ws.Range("E231").Value2 = 45383.41666666666
$ws.Range("F231").Value = "Stal Mielec"
$ws.Range("G231").Value = "Lech Poznan"
$ws.Range("K231").Value2 = 4.333
$ws.Range("L231").Value2 = 3.4
$ws.Range("M231").Value2 = 1.85
$ws.Range("N231").Value2 = 4.333
$ws.Range("O231").Value2 = 3.4
$ws.Range("P231").Value2 = 1.85
$ws.Range("Q231").Value2 = 0.5
$ws.Range("R231").Value2 = 1.95
$ws.Range("S231").Value2 = 1.9
$ws.Range("T231").Value2 = 2.25
$ws.Range("U231").Value2 = 1.925
$ws.Range("V231").Value2 = 1.925
$ws.Range("W231").Value2 = 0
$ws.Range("X231").Value2 = 0
$ws.Range("Y231").Value2 = 0
$ws.Range("Z231").Value2 = 0
$ws.Range("AA231").Value2 = 0

# Row 232 (was old row 228, match 6775578)
$ws.Range("A232").Value2 = 230
$ws.Range("B232").Value2 = 6775578
$ws.Range("C232").Value = "Poland Ekstraklasa"
$ws.Range("D232").Value = "Poland Ekstraklasa"
$ws.Range("E232").Value2 = 45383.52083333334
$ws.Range("F232").Value = "Widzew Lodz"
$ws.Range("G232").Value = "Korona Kielce"
$ws.Range("K232").Value2 = 2.25
$ws.Range("L232").Value2 = 3.2
$ws.Range("M232").Value2 = 3.4
$ws.Range("N232").Value2 = 2.2
$ws.Range("O232").Value2 = 3.25
$ws.Range("P232").Value2 = 3.4
$ws.Range("Q232").Value2 = -0.25
$ws.Range("R232").Value2 = 1.85
$ws.Range("S232").Value2 = 2
$ws.Range("T232").Value2 = 2.5
$ws.Range("U232").Value2 = 2.025
$ws.Range("V232").Value2 = 1.825
$ws.Range("W232").Value2 = 0
$ws.Range("X232").Value2 = 0
$ws.Range("Y232").Value2 = 0
$ws.Range("Z232").Value2 = 0
$ws.Range("AA232").Value2 = 0

# Row 233 (was old row 229, match 6775573)
$ws.Range("A233").Value2 = 231
$ws.Range("B233").Value2 = 6775573
$ws.Range("C233").Value = "Poland Ekstraklasa"
$ws.Range("D233").Value = "Poland Ekstraklasa"
$ws.Range("E233").Value2 = 45383.625
$ws.Range("F233").Value = "Gornik Zabrze"
$ws.Range("G233").Value = "Legia Warsaw"
$ws.Range("K233").Value2 = 3.6
$ws.Range("L233").Value2 = 3.5
$ws.Range("M233").Value2 = 2.05
$ws.Range("N233").Value2 = 3.8
$ws.Range("O233").Value2 = 3.5
$ws.Range("P233").Value2 = 1.95
$ws.Range("Q233").Value2 = 0.5
$ws.Range("R233").Value2 = 1.825
$ws.Range("S233").Value2 = 2.025
$ws.Range("T233").Value2 = 2.5
$ws.Range("U233").Value2 = 1.925
$ws.Range("V233").Value2 = 1.925
$ws.Range("W233").Value2 = 0
$ws.Range("X233").Value2 = 0
$ws.Range("Y233").Value2 = 0
$ws.Range("Z233").Value2 = 0
$ws.Range("AA233").Value2 = 0

# Row 234 (was old row 230, match 6775577)
$ws.Range("A234").Value2 = 232
$ws.Range("B234").Value2 = 6775577
$ws.Range("C234").Value = "Poland Ekstraklasa"
$ws.Range("D234").Value = "Poland Ekstraklasa"
$ws.Range("E234").Value2 = 45384.58333333334
$ws.Range("F234").Value = "Warta Poznan"
$ws.Range("G234").Value = "Zaglebie Lubin"
$ws.Range("K234").Value2 = 2.9
$ws.Range("L234").Value2 = 3.1
$ws.Range("M234").Value2 = 2.55
$ws.Range("N234").Value2 = 3
$ws.Range("O234").Value2 = 3.1
$ws.Range("P234").Value2 = 2.55
$ws.Range("Q234").Value2 = 0
$ws.Range("R234").Value2 = 2.05
$ws.Range("S234").Value2 = 1.8
$ws.Range("T234").Value2 = 2.25
$ws.Range("U234").Value2 = 2.025
$ws.Range("V234").Value2 = 1.825
$ws.Range("W234").Value2 = 0
$ws.Range("X234").Value2 = 0
$ws.Range("Y234").Value2 = 0
$ws.Range("Z234").Value2 = 0
$ws.Range("AA234").Value2 = 0
